$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Step 1: Move the "Experience" section (Heading2 "Experience" plus its
# four Heading3 sub-entries and all their bullet paragraphs -- 17
# paragraphs total, currently the very first 17 paragraphs in the body)
# so that it sits after the Education section's "Related coursework..."
# paragraph and before the "Projects" heading.
#
# The engine loses per-paragraph formatting (pStyle/numPr/rPr) when a
# Cut/Paste spans more than one paragraph at once, so the block is
# relocated one paragraph at a time. Paragraphs are cut starting from
# the last one in the block and pasted right after the anchor paragraph;
# because each new paragraph is inserted immediately after the same
# anchor, working back-to-front reproduces the original internal order.
# -----------------------------------------------------------------------

$blockCount = 17

for ($k = $blockCount; $k -ge 1; $k--) {
    $cutPara = $d.Paragraphs.Item($k)
    $cutPara.Range.Cut()

    $anchorIndex = -1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -like "Related coursework*") {
            $anchorIndex = $i
            break
        }
    }

    $targetPara = $d.Paragraphs.Item($anchorIndex)
    $ins = $d.Range($targetPara.Range.End, $targetPara.Range.End)
    $ins.Paste()
}

# -----------------------------------------------------------------------
# Step 2: Education bullets.
# -----------------------------------------------------------------------

# "Major: Computer Systems Engineering " + "with a GPA of 3.6." -> merged,
# GPA bumped to 3.7.
$d.Content.Find.Execute(
    "Major: Computer Systems Engineering with a GPA of 3.6.", $false, $false,
    $false, $false, $false, $true, 1, $false,
    "Major: Computer Systems Engineering with a GPA of 3.7.", 2) | Out-Null

# "Related coursework: ..." runs consolidate into a single run (same text).
$d.Content.Find.Execute(
    "Related coursework: Computer architecture, Object-oriented development, Data structures and algorithms, Intro to machine learning, Network and cyber security, Communications engineering, Project management",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Related coursework: Computer architecture, Object-oriented development, Data structures and algorithms, Intro to machine learning, Network and cyber security, Communications engineering, Project management",
    2) | Out-Null

# -----------------------------------------------------------------------
# Step 3: Skills section.
# -----------------------------------------------------------------------

# Languages: append the new language/platform entries.
$d.Content.Find.Execute(
    "Python, Java, C/C++, Verilog", $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Python, Java, C/C++, Verilog, Linux/Unix Kernel, Windows Command Line, Bash",
    2) | Out-Null

# Locate the Skills bullet paragraphs (NoSpacing + numId 30) by index so
# the edits are unambiguous even though several bullets start similarly.
$skillsParas = @()
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text
    if ($txt -like "Operating Systems:*" -or $txt -like "Tools:*" -or $txt -like "Technologies:*") {
        $skillsParas += $i
    }
}

# skillsParas[0] = "Operating Systems: Linux, Windows"       -> becomes Tools
# skillsParas[1] = "Tools: Git, JIRA, TeamCity, BitBucket"   -> becomes Leadership
# skillsParas[2] = "Technologies: Docker, ..."               -> removed

$osParaIndex = $skillsParas[0]
$osPara = $d.Paragraphs.Item($osParaIndex)
$osPara.Range.Text = "Tools: Git, JIRA, TeamCity, BitBucket, Docker, Maven, PyTest, OpenSSL, Embedded / SoC Dev, Machine Learning (CNN), Computer Vision, TensorFlow, OpenCV`r"
$osPara.Range.Font.Color = -587137025
$labelEnd = $osPara.Range.Start + 5
$labelRange = $d.Range($osPara.Range.Start, $labelEnd)
$labelRange.Style = "Heading5Char"

$toolsParaIndex = $skillsParas[1]
$toolsPara = $d.Paragraphs.Item($toolsParaIndex)
$toolsPara.Range.Text = "Leadership: Skills gained as an Air Cadet Warrant Officer and as a Teaching Assistant, mentored students to facilitate learning `r"
$labelEnd2 = $toolsPara.Range.Start + 10
$labelRange2 = $d.Range($toolsPara.Range.Start, $labelEnd2)
$labelRange2.Style = "Heading5Char"

$technologiesParaIndex = $skillsParas[2]
$technologiesPara = $d.Paragraphs.Item($technologiesParaIndex)
$technologiesPara.Range.Delete()

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output ("[$i] " + $d.Paragraphs.Item($i).Range.Text)
}
